$d = $word.ActiveDocument

# --- Change 1: first bullet under "2D or 3D?" ---------------------------
$d.Content.Find.Execute(
    "Depending on amount of data", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "First of all, 2D (3D will be done if we have time at the end)", 2) | Out-Null

# --- Change 2: "Mov" + bookmark + "ing nodes" -> single "Moving nodes" run,
#     bookmark relocated to the very end of the list (after "D3.js:"),
#     and the trailing empty list paragraph removed. ----------------------

# 2a. Drop the existing (hidden) _GoBack bookmark so the split runs can be
#     merged back into a single run.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2b. Force Word to re-merge the adjacent "Mov"/"ing nodes" runs into one
#     contiguous run by replacing the (now bookmark-free) text with itself.
$d.Content.Find.Execute(
    "Moving nodes", $true, $false, $false, $false, $false,
    $true, 1, $false, "Moving nodes", 2) | Out-Null

# 2c. Remove the now-empty trailing list paragraph that used to hold the
#     bookmark's paragraph (directly after the "D3.js:" paragraph).
$paras = @($d.Paragraphs)
for ($i = 0; $i -lt $paras.Count; $i++) {
    if ($paras[$i].Range.Text.TrimEnd([char]13) -eq "D3.js:") {
        $paras[$i + 1].Range.Delete()
        break
    }
}

# 2d. Re-insert the _GoBack bookmark right after "D3.js:" (collapsed, before
#     the paragraph mark). A bookmark can't be created directly abutting a
#     paragraph mark, so a temporary placeholder character is used to hold
#     the spot, the bookmark is anchored next to it, then the placeholder is
#     removed - the (collapsed) bookmark stays put.
$paras = @($d.Paragraphs)
for ($i = 0; $i -lt $paras.Count; $i++) {
    if ($paras[$i].Range.Text.TrimEnd([char]13) -eq "D3.js:") {
        $p = $paras[$i]
        $endPos = $p.Range.End - 1
        $d.Range($endPos, $endPos).InsertAfter("X")

        $d.Bookmarks.Add("_GoBack", $d.Range($endPos, $endPos))

        $d.Range($endPos, $endPos + 1).Delete()
        break
    }
}
